# Append the four new leetcode-id rows to column A, right after the
# existing data (which currently ends at row 104), then move the
# selection to the new last cell (A108) so the view reflects where the
# user was last working - matching Excel's normal "type value, Enter"
# workflow down a single column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A105").Value = 435
$ws.Range("A106").Value = 452
$ws.Range("A107").Value = 1024
$ws.Range("A108").Value = 64

$ws.Range("A108").Select()
